# Replace the figure caption/placeholder paragraph near the top of the
# document ("I cut him short with: "Take that, you monster!"-and, jumping
# back, I press the switch. There is a blinding flash.](images/killing_flash.png)")
# with a simple "insert figure here" placeholder, as requested in the
# commit "new word docs with figure numbers".

$d = $word.ActiveDocument

# Find the specific paragraph: it is the only one whose text ends with the
# markdown-style image reference "...killing_flash.png)". (A later
# paragraph in the story also starts with "I cut him short with:", so we
# must be careful to only touch this one.)
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "I cut him short with:*killing_flash.png)*") {
        $target = $p.Range
        break
    }
}

if ($null -eq $target) {
    throw "Could not locate the killing_flash.png caption paragraph"
}

$target.Find.ClearFormatting()
$target.Find.Execute(
    "I cut him short with: “Take that, you monster!”—and, jumping back, I press the switch. There is a blinding flash.](images/killing_flash.png)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[INSERT FIGURE 66.1 NEAR HERE]", 2)
